# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 17 de Mayo de 2020 a las 23:35"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1524463
$ws.Range("C4").Value = 16690
$ws.Range("D4").Value = 343150
$ws.Range("E4").Value = 1090412
$ws.Range("G4").Value = 788
$ws.Range("H4").Value = 90901

# Row 7 - Reino Unido
$ws.Range("B7").Value = 243695
$ws.Range("C7").Value = 3534

# Row 8 - Brasil
$ws.Range("B8").Value = 236131
$ws.Range("C8").Value = 2989
$ws.Range("E8").Value = 130683
$ws.Range("G8").Value = 143
$ws.Range("H8").Value = 15776

# Row 11 - Alemania
$ws.Range("B11").Value = 176625
$ws.Range("C11").Value = 381
$ws.Range("E11").Value = 15177
$ws.Range("G11").Value = 21
$ws.Range("H11").Value = 8048

# Row 17 - Canada
$ws.Range("B17").Value = 77002
$ws.Range("C17").Value = 1138
$ws.Range("D17").Value = 38550
$ws.Range("E17").Value = 32670
$ws.Range("G17").Value = 103
$ws.Range("H17").Value = 5782
